$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.032238159739954
$ws.Range("D2").Value = 1.035411514214439
$ws.Range("E2").Value = 1.040941993116917
$ws.Range("F2").Value = 1.050074940190408
$ws.Range("I2").Value = 1.033533150488972
$ws.Range("J2").Value = 1.037369480076597
$ws.Range("K2").Value = 1.038208195315353
$ws.Range("L2").Value = 1.043722913166294
$ws.Range("M2").Value = 1.052830226853501
$ws.Range("N2").Value = 1.005712725503983
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.033222425333232
$ws.Range("D3").Value = 1.036134189439881
$ws.Range("E3").Value = 1.041823143574659
$ws.Range("F3").Value = 1.051055720509156
$ws.Range("I3").Value = 1.033693886621651
$ws.Range("J3").Value = 1.03799577980256
$ws.Range("K3").Value = 1.038740490496891
$ws.Range("L3").Value = 1.044414404366428
$ws.Range("M3").Value = 1.05362294654327
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.03385977406638
$ws.Range("D4").Value = 1.036601951805587
$ws.Range("E4").Value = 1.042394089855593
$ws.Range("F4").Value = 1.051691228812076
$ws.Range("I4").Value = 1.033796635727209
$ws.Range("J4").Value = 1.038400907677486
$ws.Range("K4").Value = 1.039084395181607
$ws.Range("L4").Value = 1.044861991377187
$ws.Range("M4").Value = 1.054136149083095
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.034127825593162
$ws.Range("D5").Value = 1.036798632086645
$ws.Range("E5").Value = 1.042634301950496
$ws.Range("F5").Value = 1.051958605483188
$ws.Range("I5").Value = 1.03383952980405
$ws.Range("J5").Value = 1.038571191472868
$ws.Range("K5").Value = 1.039228845869806
$ws.Range("L5").Value = 1.045050190956526
$ws.Range("M5").Value = 1.054351960606277
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.034172839059974
$ws.Range("D6").Value = 1.036831657467779
$ws.Range("E6").Value = 1.042674645507492
$ws.Range("F6").Value = 1.052003511423274
$ws.Range("I6").Value = 1.033846714205096
$ws.Range("J6").Value = 1.038599780972136
$ws.Range("K6").Value = 1.039253092310568
$ws.Range("L6").Value = 1.045081792460368
$ws.Range("M6").Value = 1.054388199867704
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.033863355352227
$ws.Range("D7").Value = 1.036604579729157
$ws.Range("E7").Value = 1.04239729885089
$ws.Range("F7").Value = 1.051694800692567
$ws.Range("I7").Value = 1.033797210065852
$ws.Range("J7").Value = 1.038403183144563
$ws.Range("K7").Value = 1.039086325838355
$ws.Range("L7").Value = 1.044864505976764
$ws.Range("M7").Value = 1.05413903252816
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.032570700574358
$ws.Range("D8").Value = 1.035655715268952
$ws.Range("E8").Value = 1.04123961920369
$ws.Range("F8").Value = 1.05040621701599
$ws.Range("I8").Value = 1.033587732178369
$ws.Range("J8").Value = 1.037581167387216
$ws.Range("K8").Value = 1.038388195287138
$ws.Range("L8").Value = 1.043956575063674
$ws.Range("M8").Value = 1.053098075743804
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.030296457052602
$ws.Range("D9").Value = 1.033984848373055
$ws.Range("E9").Value = 1.03920568861288
$ws.Range("F9").Value = 1.048142343295738
$ws.Range("I9").Value = 1.033208993213326
$ws.Range("J9").Value = 1.036131713349781
$ws.Range("K9").Value = 1.037154013877927
$ws.Range("L9").Value = 1.042357844361874
$ws.Range("M9").Value = 1.051265813390655
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.028782743725126
$ws.Range("D10").Value = 1.032871788747539
$ws.Range("E10").Value = 1.037853870803772
$ws.Range("F10").Value = 1.046637718527405
$ws.Range("I10").Value = 1.032950067206241
$ws.Range("J10").Value = 1.035164817161203
$ws.Range("K10").Value = 1.036328601411076
$ws.Range("L10").Value = 1.041292858766433
$ws.Range("M10").Value = 1.050045735426714
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.028127878424784
$ws.Range("D11").Value = 1.032390039436615
$ws.Range("E11").Value = 1.037269513991818
$ws.Range("F11").Value = 1.045987310412934
$ws.Range("I11").Value = 1.032836429526945
$ws.Range("J11").Value = 1.034746009930003
$ws.Range("K11").Value = 1.035970577510792
$ws.Range("L11").Value = 1.040831917618992
$ws.Range("M11").Value = 1.049517779979842
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.02788472042617
$ws.Range("D12").Value = 1.032211129514218
$ws.Range("E12").Value = 1.03705260758633
$ws.Range("F12").Value = 1.045745886895774
$ws.Range("I12").Value = 1.032793991374343
$ws.Range("J12").Value = 1.034590426810187
$ws.Range("K12").Value = 1.035837499997513
$ws.Range("L12").Value = 1.040660735175555
$ws.Range("M12").Value = 1.049321726825227
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.02793687464431
$ws.Range("D13").Value = 1.032249504776605
$ws.Range("E13").Value = 1.037099127952987
$ws.Range("F13").Value = 1.045797665484201
$ws.Range("I13").Value = 1.032803104814947
$ws.Range("J13").Value = 1.03462380078581
$ws.Range("K13").Value = 1.035866049702102
$ws.Range("L13").Value = 1.040697452945705
$ws.Range("M13").Value = 1.049363778478226
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.028107777107361
$ws.Range("D14").Value = 1.032375250001724
$ws.Range("E14").Value = 1.037251581379723
$ws.Range("F14").Value = 1.045967350856238
$ws.Range("I14").Value = 1.032832926225702
$ws.Range("J14").Value = 1.034733149773351
$ws.Range("K14").Value = 1.035959579143427
$ws.Range("L14").Value = 1.040817766974659
$ws.Range("M14").Value = 1.049501573078016
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.028213087451969
$ws.Range("D15").Value = 1.032452730224379
$ws.Range("E15").Value = 1.037345532841941
$ws.Range("F15").Value = 1.046071921778504
$ws.Range("I15").Value = 1.03285126997475
$ws.Range("J15").Value = 1.034800520737539
$ws.Range("K15").Value = 1.036017193621334
$ws.Range("L15").Value = 1.040891900635253
$ws.Range("M15").Value = 1.049586479931501
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.028826217304064
$ws.Range("D16").Value = 1.032903765453369
$ws.Range("E16").Value = 1.037892673577627
$ws.Range("F16").Value = 1.046680907333017
$ws.Range("I16").Value = 1.032957576957702
$ws.Range("J16").Value = 1.035192609244746
$ws.Range("K16").Value = 1.036352349398967
$ws.Range("L16").Value = 1.041323454298723
$ws.Range("M16").Value = 1.050080781471029
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.029210973857527
$ws.Range("D17").Value = 1.033186745800674
$ws.Range("E17").Value = 1.038236146134155
$ws.Range("F17").Value = 1.047063204322041
$ws.Range("I17").Value = 1.033023853623362
$ws.Range("J17").Value = 1.035438520470568
$ws.Range("K17").Value = 1.036562419770016
$ws.Range("L17").Value = 1.041594211917702
$ws.Range("M17").Value = 1.05039093718177
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.029435452013018
$ws.Range("D18").Value = 1.033351823903438
$ws.Range("E18").Value = 1.038436583235917
$ws.Range("F18").Value = 1.047286298166507
$ws.Range("I18").Value = 1.033062364895124
$ws.Range("J18").Value = 1.035581943268992
$ws.Range("K18").Value = 1.036684890879717
$ws.Range("L18").Value = 1.041752159942221
$ws.Range("M18").Value = 1.05057187911435
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.029512002771529
$ws.Range("D19").Value = 1.033408114701976
$ws.Range("E19").Value = 1.038504943250503
$ws.Range("F19").Value = 1.047362385429373
$ws.Range("I19").Value = 1.033075471332279
$ws.Range("J19").Value = 1.035630844497951
$ws.Range("K19").Value = 1.036726640263953
$ws.Range("L19").Value = 1.041806019449312
$ws.Range("M19").Value = 1.050633581265374
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.029169687306525
$ws.Range("D20").Value = 1.033156382562867
$ws.Range("E20").Value = 1.038199284906722
$ws.Range("F20").Value = 1.047022176451553
$ws.Range("I20").Value = 1.033016757943836
$ws.Range("J20").Value = 1.035412137868376
$ws.Range("K20").Value = 1.036539887343237
$ws.Range("L20").Value = 1.041565160145006
$ws.Range("M20").Value = 1.050357656952307
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.028057448153838
$ws.Range("D21").Value = 1.032338220239989
$ws.Range("E21").Value = 1.037206683503038
$ws.Range("F21").Value = 1.045917378124154
$ws.Range("I21").Value = 1.032824150855386
$ws.Range("J21").Value = 1.034700949753033
$ws.Range("K21").Value = 1.035932039572782
$ws.Range("L21").Value = 1.040782336611906
$ws.Range("M21").Value = 1.049460994531817
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.027358649462838
$ws.Range("D22").Value = 1.031824001939178
$ws.Range("E22").Value = 1.036583462465288
$ws.Range("F22").Value = 1.04522371441645
$ws.Range("I22").Value = 1.032701731883061
$ws.Range("J22").Value = 1.034253685319358
$ws.Range("K22").Value = 1.035549332623949
$ws.Range("L22").Value = 1.040290327777487
$ws.Range("M22").Value = 1.04889753443449
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.02772904730401
$ws.Range("D23").Value = 1.032096580097478
$ws.Range("E23").Value = 1.036913761188507
$ws.Range("F23").Value = 1.045591346657797
$ws.Range("I23").Value = 1.032766753393796
$ws.Range("J23").Value = 1.034490799034381
$ws.Range("K23").Value = 1.035752262661276
$ws.Range("L23").Value = 1.040551133293232
$ws.Range("M23").Value = 1.049196205854446
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.029188342748205
$ws.Range("D24").Value = 1.033170102339977
$ws.Range("E24").Value = 1.038215940614136
$ws.Range("F24").Value = 1.047040714851278
$ws.Range("I24").Value = 1.033019964629561
$ws.Range("J24").Value = 1.035424059071078
$ws.Range("K24").Value = 1.036550068960551
$ws.Range("L24").Value = 1.041578287329801
$ws.Range("M24").Value = 1.050372694752477
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.030883974986472
$ws.Range("D25").Value = 1.034416662974904
$ws.Range("E25").Value = 1.039730784495485
$ws.Range("F25").Value = 1.048726798899576
$ws.Range("I25").Value = 1.033308042388879
$ws.Range("J25").Value = 1.036506540242276
$ws.Range("K25").Value = 1.037473546049255
$ws.Range("L25").Value = 1.042771011466524
$ws.Range("M25").Value = 1.051739249508672
